$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rewrite the data grid (A1:C34) to the final, sorted layout ---
$ws.Range('A1').Value = 'hebrew_text'
$ws.Range('B1').Value = 'hebrew_option'
$ws.Range('C1').Value = 'label'

$ws.Range('A2').Value = 'מאז הסקר האחרון, דעתו/ה של הילד/ה שלי הוסחה בקלות'
$ws.Range('B2').Value = ''
$ws.Range('C2').Value = 'P_ADHD_Distracted'

$ws.Range('A3').Value = 'למיטב ידיעתי - ברגע זה, הילד/ה שלי מרגיש/ה חסר/ת מנוחה'
$ws.Range('B3').Value = ''
$ws.Range('C3').Value = 'P_ADHD_Restless'

$ws.Range('A4').Value = 'מאז הסקר האחרון, הילד/ה שלי התפרץ/ה בכעס באחת או יותר מהדרכים הבאות'
$ws.Range('B4').Value = 'הרביץ/ה או פגע/ה במישהו'
$ws.Range('C4').Value = 'P_Agr_hit'

$ws.Range('A5').Value = 'מאז הסקר האחרון, הילד/ה שלי התפרץ/ה בכעס באחת או יותר מהדרכים הבאות'
$ws.Range('B5').Value = 'לא התפרץ/ה בכעס'
$ws.Range('C5').Value = 'P_Agr_none'

$ws.Range('A6').Value = 'מאז הסקר האחרון, הילד/ה שלי כעס/ה או התעצבן/ה כשדברים לא קרו כמו שהוא/היא רצה/תה'
$ws.Range('B6').Value = ''
$ws.Range('C6').Value = 'P_Agr_NotAsWant'

$ws.Range('A7').Value = 'מאז הסקר האחרון, הילד/ה שלי התפרץ/ה בכעס באחת או יותר מהדרכים הבאות'
$ws.Range('B7').Value = 'אחר (בבקשה לפרט)'
$ws.Range('C7').Value = 'P_Agr_other'

$ws.Range('A8').Value = 'מאז הסקר האחרון, הילד/ה שלי התפרץ/ה בכעס באחת או יותר מהדרכים הבאות'
$ws.Range('B8').Value = 'טרק/ה דלת'
$ws.Range('C8').Value = 'P_Agr_slam'

$ws.Range('A9').Value = 'מאז הסקר האחרון, הילד/ה שלי התפרץ/ה בכעס באחת או יותר מהדרכים הבאות'
$ws.Range('B9').Value = 'זרק/ה משהו'
$ws.Range('C9').Value = 'P_Agr_throw_smt'

$ws.Range('A10').Value = 'מאז הסקר האחרון, הילד/ה שלי התפרץ/ה בכעס באחת או יותר מהדרכים הבאות'
$ws.Range('B10').Value = 'זרק/ה משהו על מישהו'
$ws.Range('C10').Value = 'P_Agr_throw_twd'

$ws.Range('A11').Value = 'מאז הסקר האחרון, הילד/ה שלי התפרץ/ה בכעס באחת או יותר מהדרכים הבאות'
$ws.Range('B11').Value = 'צעק/ה'
$ws.Range('C11').Value = 'P_Agr_yelled'

$ws.Range('A12').Value = 'למיטב ידיעתי - ברגע זה, הילד/ה שלי מרגיש/ה מרוגז/ת או כועס/ת'
$ws.Range('B12').Value = ''
$ws.Range('C12').Value = 'P_Angry_now'

$ws.Range('A13').Value = 'למיטב ידיעתי - ברגע זה, הילד/ה שלי מרגיש/ה פחד או לחץ'
$ws.Range('B13').Value = ''
$ws.Range('C13').Value = 'P_Anx_now'

$ws.Range('A14').Value = 'מאז הסקר האחרון, הילד/ה שלי הרגיש/ה מודאג/ת או מפוחד/ת'
$ws.Range('B14').Value = ''
$ws.Range('C14').Value = 'P_Anx_Worry'

$ws.Range('A15').Value = 'מאז הסקר האחרון, לילד/ה שלי היה קשה להפסיק לעשות משהו כאשר ביקשו ממנו/ה להפסיק'
$ws.Range('B15').Value = ''
$ws.Range('C15').Value = 'P_Difficult_Stop'

$ws.Range('A16').Value = 'מאז הסקר האחרון איימתי להעניש את הילד/ה שלי ואז לא עשיתי זאת'
$ws.Range('B16').Value = 'לא'
$ws.Range('C16').Value = 'P_Discipline'

$ws.Range('A17').Value = 'מאז הסקר האחרון איימתי להעניש את הילד/ה שלי ואז לא עשיתי זאת'
$ws.Range('B17').Value = 'כן'
$ws.Range('C17').Value = 'P_Discipline'

$ws.Range('A18').Value = 'מאז הסקר האחרון, יצא שהילד/ה שלי אמר/ה את הדבר הראשון שעלה לו/ה לראש מבלי לעצור ולחשוב'
$ws.Range('B18').Value = ''
$ws.Range('C18').Value = 'P_IC_FirstOnMind'

$ws.Range('A19').Value = 'היום אני והילד/ה שלי בילינו זמן כיף ביחד'
$ws.Range('B19').Value = ''
$ws.Range('C19').Value = 'P_Inv_Fun'

$ws.Range('A20').Value = 'היום עזרתי בפעילות שבה הילד/ה שלי מעורב, או עזרתי לו/ה במשהו כללי אחר (כמו הקפצה לחוג, עזרה בשיעורי בית וכו'')'
$ws.Range('B20').Value = ''
$ws.Range('C20').Value = 'P_Inv_Help'

$ws.Range('A21').Value = 'היום דיברתי עם הילד/ה שלי על היום שלו/ה או על החברים שלו/ה'
$ws.Range('B21').Value = ''
$ws.Range('C21').Value = 'P_Inv_Talk'

$ws.Range('A22').Value = 'מאז הסקר האחרון, הילד/ה שלי הרגיש/ה מתוסכל/ת'
$ws.Range('B22').Value = ''
$ws.Range('C22').Value = 'P_Irr_Frustration'

$ws.Range('A23').Value = 'למיטב ידיעתי - ברגע זה הילד/ה שלי מרגיש/ה טוב'
$ws.Range('B23').Value = ''
$ws.Range('C23').Value = 'P_Mood_Good'

$ws.Range('A24').Value = 'מאז הסקר האחרון, הילד/ה שלי הרגיש/ה עצוב/ה או מדוכא/ת'
$ws.Range('B24').Value = ''
$ws.Range('C24').Value = 'P_Mood_Sad'

$ws.Range('A25').Value = 'מאז הסקר האחרון, הילד/ה שלי עצבן/ה אותי'
$ws.Range('B25').Value = ''
$ws.Range('C25').Value = 'P_PC_Annoy'

$ws.Range('A26').Value = 'מאז הסקר האחרון, העברתי ביקורת על הילד/ה שלי'
$ws.Range('B26').Value = ''
$ws.Range('C26').Value = 'P_PC_Criticism'

$ws.Range('A27').Value = 'מאז הסקר האחרון, הילד/ה שלי שיתף/ה אותי ברגשות/תחושות שלו/ה'
$ws.Range('B27').Value = ''
$ws.Range('C27').Value = 'P_PC_Sharing'

$ws.Range('A28').Value = 'מאז הסקר האחרון, החמאתי לילד/ה שלי כשהוא/היא עשה/תה משהו בצורה טובה'
$ws.Range('B28').Value = ''
$ws.Range('C28').Value = 'P_Positive'

$ws.Range('A29').Value = 'מאז הסקר האחרון, הסכמתי לילד/ה שלי כל מה שרצ/תה'
$ws.Range('B29').Value = ''
$ws.Range('C29').Value = 'P_PS_Agree'

$ws.Range('A30').Value = 'מאז הסקר האחרון, הרגשתי בטוח/ה ביכולותיי ההוריות'
$ws.Range('B30').Value = ''
$ws.Range('C30').Value = 'P_PS_Confident'

$ws.Range('A31').Value = 'מאז הסקר האחרון, התעצבנתי או צעקתי על הילד/ה שלי'
$ws.Range('B31').Value = ''
$ws.Range('C31').Value = 'P_PS_GotAngry'

$ws.Range('A32').Value = 'מאז הסקר האחרון, הייתי סבלני/ת כלפי הילד/ה שלי'
$ws.Range('B32').Value = ''
$ws.Range('C32').Value = 'P_PS_Patient'

$ws.Range('A33').Value = 'תן/י דוגמאות למשהו שעצבן את הילד/ה שלך היום (טריגר) (אפשר גם להקליט תשובה)'
$ws.Range('B33').Value = ''
$ws.Range('C33').Value = 'P_triggres'

$ws.Range('A34').Value = 'תן/י דוגמאות למשהו שעצבן את הילד/ה שלך היום (טריגר)'
$ws.Range('B34').Value = ''
$ws.Range('C34').Value = 'P_triggres'
# --- Header row formatting (bold) ---
$ws.Range('A1:C1').Font.Bold = $true

# --- Turn the range into a proper Excel Table ---
$tbl = $ws.ListObjects.Add(1, $ws.Range('A1:C34'), $null, 1)
$tbl.Name = 'Table1'

# --- Keep the table sorted by the "label" column (data already in this order) ---
$tbl.Sort.SortFields.Clear()
[void]$tbl.Sort.SortFields.Add($ws.Range('C2:C34'))
$tbl.Sort.Header = 0
[void]$tbl.Sort.Apply()

# --- Reset the active selection back to A1 ---
[void]$ws.Range('A1').Select()
